$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 695; everything from old row 695 downward
# shifts down by one (old 695 -> new 696, ..., old 772 -> new 773).
$ws.Rows.Item(695).Insert()

# Populate the newly inserted row 695 with its data.
$ws.Cells.Item(695, 1).Value = 6
$ws.Cells.Item(695, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(695, 3).Value = "Metropolitana"
$ws.Cells.Item(695, 4).Value = 44776
$ws.Cells.Item(695, 5).Value = 13
$ws.Cells.Item(695, 6).Value = 100112003
$ws.Cells.Item(695, 7).Value = "Ajo"
$ws.Cells.Item(695, 8).Value = "Chino"
$ws.Cells.Item(695, 9).Value = "Primera"
$ws.Cells.Item(695, 10).Value = 3700
$ws.Cells.Item(695, 11).Value = 23000
$ws.Cells.Item(695, 12).Value = 24000
$ws.Cells.Item(695, 13).Value = 23595
$ws.Cells.Item(695, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(695, 15).Value = "China"
$ws.Cells.Item(695, 16).Value = 2360
$ws.Cells.Item(695, 17).Value = 10
$ws.Cells.Item(695, 18).Value = "Hortaliza"
